$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 4 de Agosto de 2020 a las 23:06"

# Helper to set a whole data row (columns A-H) in one shot
function Set-Row {
    param($row, $country, $total, $nuevos, $activos, $recuperados, $criticos, $muertesHoy, $muertes)
    $ws.Range("A$row").Value = $country
    $ws.Range("B$row").Value = $total
    $ws.Range("C$row").Value = $nuevos
    $ws.Range("D$row").Value = $activos
    $ws.Range("E$row").Value = $recuperados
    $ws.Range("F$row").Value = $criticos
    $ws.Range("G$row").Value = $muertesHoy
    $ws.Range("H$row").Value = $muertes
}

# Rows whose country/label stays the same but the statistics were refreshed
Set-Row 4   "Estados Unidos"          4906672 42756 2474672 2272030 0 1042 159970
Set-Row 5   "Brasil"                  2801921 50256 1912319  793783 0 1117  95819
Set-Row 8   "Sudafrica"                521318  4456  363751  148683 0  345   8884
Set-Row 36  "Israel"                    76198  1768   49834   25803 0   15    561
Set-Row 71  "Costa Rica"                19837   435    6590   13066 0   10    181
Set-Row 76  "Costa de Marfil"           16293    73   11955    4235 0    1    103
Set-Row 80  "Estado de Palestina"       12770   229    6419    6265 0    2     86
Set-Row 97  "Mauritania"                 6418    36    5209    1052 0    0    157
Set-Row 105 "Republica de Africa Central" 4618    4    1640    2919 0    0     59
Set-Row 109 "Zimbabue"                   4221   146    1238    2902 0    1     81
Set-Row 128 "Ruanda"                     2099     7    1222     872 0    0      5
Set-Row 146 "Republica de Chipre"        1180    25     856     305 0    0     19
Set-Row 168 "Comoras"                     388     2     340      41 0    0      7
Set-Row 185 "Aruba"                       132     8     112      17 0    0      3

# Rows where two countries swapped rank order because of the refreshed totals
Set-Row 93  "Guinea"     7489 125 6591  850 0 2  48
Set-Row 94  "Finlandia"  7483  17 6950  202 0 2 331

Set-Row 139 "Angola"   1344 64  503 782 0 1 59
Set-Row 140 "Uruguay"  1291  0 1023 232 0 0 36

Set-Row 161 "Gambia"    671  0   79 578 0 5 14
Set-Row 162 "Vietnam"   670 18  378 284 0 2  8
Set-Row 163 "Reunion"   667  0  592  71 0 0  4
Set-Row 164 "Tanzania"  509  0  183 305 0 0 21
